$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "InteresRefinanciar" column (was column J). Excel shifts
#    every column after it one slot to the left (K->J, L->K, ... W->V) and
#    the used range shrinks from A1:W7 to A1:V7.
# ---------------------------------------------------------------------------
$ws.Range("J1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Update the data row (row 2) with the new proposal's values.
#    (Columns below use the NEW, post-deletion lettering.)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "1940821"
$ws.Range("H2").Value = "080-01-9434661"
$ws.Range("I2").Value = "808.46"
$ws.Range("N2").Value = "6"
$ws.Range("O2").Value = "42.58"

# U2 (previously the empty V2) now receives a value. Borrow the text-number
# format from A2 first so the new value is stored as text (matching how the
# rest of the "numeric looking" fields in this row are stored), then drop
# back to the default "Normal" style so no explicit style index is written,
# just like the sibling empty cells around it.
$ws.Range("A2").Copy()
$ws.Range("U2").PasteSpecial(-4122)
$ws.Range("U2").Value = "4899906"
$ws.Range("U2").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Formatting tweaks that came along with the new data.
# ---------------------------------------------------------------------------
# F2 ("SIN PROMOCION") is highlighted with the green fill used elsewhere in
# the sheet.
$ws.Range("F2").Interior.Color = 5296274

# H2 ("Pagare" value) loses its previous yellow highlight - copy the plain
# text-number format from I2 so it matches the rest of the row.
$ws.Range("I2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# Row 3 (the blank placeholder row) gets center-aligned cells from A3 to T3.
$ws.Range("A3:T3").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Leave the selection on the last column of the new layout.
# ---------------------------------------------------------------------------
$ws.Range("V2").Select()
